$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------------
# Shape "Text Box 30" (id 24) - "Space-Time Guess" -> "Initial Space-Time Guess"
# ---------------------------------------------------------------------------
$shpGuess = $s.Shapes.Item(4)
$shpGuess.Left   = 5.281685039370079
$shpGuess.Top    = 161.6840472440945
$shpGuess.Width  = 208.44790551181103
$shpGuess.Height = 31.50475590551181
$shpGuess.Fill.Visible = $false
$shpGuess.TextFrame.MarginLeft = 0
$shpGuess.TextFrame.MarginRight = 0
$shpGuess.TextFrame.TextRange.Text = "Initial Space-Time Guess"
$shpGuess.TextFrame.TextRange.Font.Size = 20
$shpGuess.TextFrame.TextRange.Font.Bold = $true

# ---------------------------------------------------------------------------
# Shape "Text Box 30" (id 52) - "Iterate"
# ---------------------------------------------------------------------------
$shpIterate = $s.Shapes.Item(10)
$shpIterate.Left   = 247.9556220472441
$shpIterate.Top    = 161.65270866141734
$shpIterate.Width  = 211.05656692913388
$shpIterate.Height = 31.50475590551181
$shpIterate.Fill.Visible = $false
$shpIterate.TextFrame.TextRange.Text = "Iterate"
$shpIterate.TextFrame.TextRange.Font.Size = 20

# ---------------------------------------------------------------------------
# Shape "Line 25" (id 66) - arrow between "Iterate" and "Converge"
# Duplicate it first (before repositioning) so the copy keeps matching
# geometry/line formatting/text for the new arrow added near the bottom of
# the deck, then move the original to its new spot.
# ---------------------------------------------------------------------------
$shpLine66 = $s.Shapes.Item(19)
$shpLineDup = $shpLine66.Duplicate()
$shpLineDup.Left   = 458.24703937007877
$shpLineDup.Top    = 180.52152755905513
$shpLineDup.Width  = 51.75278740157481
$shpLineDup.Height = 0.5491653543307087

$shpLine66.Left   = 217.36577952755906
$shpLine66.Top    = 181.07066141732284
$shpLine66.Width  = 51.75278740157481
$shpLine66.Height = 0.5491653543307087

# ---------------------------------------------------------------------------
# Shape "Text Box 30" (id 87) - "Finish" -> "Converge"
# ---------------------------------------------------------------------------
$shpFinish = $s.Shapes.Item(23)
$shpFinish.Left   = 486.07940157480317
$shpFinish.Top    = 162.36696062992127
$shpFinish.Width  = 210.74373228346457
$shpFinish.Height = 31.50475590551181
$shpFinish.Fill.Visible = $false
$shpFinish.TextFrame.TextRange.Text = "Converge"
$shpFinish.TextFrame.TextRange.Font.Size = 20

# ---------------------------------------------------------------------------
# Remove the old "Line 25" (id 127) shape - superseded by the duplicated
# line created above.
# ---------------------------------------------------------------------------
$shpOldLine = $s.Shapes.Item(29)
$shpOldLine.Delete()
